$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.023
$ws.Range("C4").Value = -11.47669999999999

$ws.Range("C5").Value = -14.46830000000001

$ws.Range("A6").Value = -21.46880000000002

$ws.Range("A7").Value = -21.51980000000002

$ws.Range("C8").Value = -11.7314

$ws.Range("A16").Value = -20.34889999999999
$ws.Range("C16").Value = -11.8864

$ws.Range("A20").Value = -22.84280000000001

$ws.Range("C22").Value = -11.1959
